# Vibration tests for ADXL1002
# Appends a new block of "Periodic Test" rows (device = ADXL1002) below the
# existing ADXL354 block (which ends at row 32) in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Frequency(C), Voltage(D), Max(E), RMS(F), std(G)
$rows = @(
    @(53,   0.625, 28,   1,  1),
    @(95,   1.25,  27,   1,  1),
    @(553,  3.1,   56,   1,  2),
    @(595,  3.2,   61,   1,  2),
    @(1053, 6.1,   75,   2,  2),
    @(1095, 6.2,   96,   2,  2),
    @(1553, 10.5,  208,  5,  5),
    @(1595, 10.6,  241,  6,  6),
    @(2053, 20,    511,  13, 13),
    @(2095, 20,    299,  7,  7),
    @(3053, 20,    714,  22, 21),
    @(3095, 20,    1078, 24, 24),
    @(4053, 20,    135,  3,  4),
    @(4095, 20,    115,  3,  3)
)

$startRow = 33
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 2).Value = "ADXL1002"
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}

# Match the author's final selection/scroll position recorded in the saved
# workbook (view state only, no data impact).
$ws.Range("G47").Select() | Out-Null
